$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; existing rows 31..57 shift down to 32..58
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly record
$ws.Cells.Item(31, 1).Value = 5
$ws.Cells.Item(31, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(31, 3).Value = "Maule"
$ws.Cells.Item(31, 4).Value = 45264
$ws.Cells.Item(31, 5).Value = 7
$ws.Cells.Item(31, 6).Value = "Fruta"
$ws.Cells.Item(31, 7).Value = 100103
$ws.Cells.Item(31, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(31, 9).Value = 100103003
$ws.Cells.Item(31, 10).Value = "Damasco"
$ws.Cells.Item(31, 11).Value = "Dina"
$ws.Cells.Item(31, 12).Value = "Segunda"
$ws.Cells.Item(31, 13).Value = 150
$ws.Cells.Item(31, 14).Value = 17000
$ws.Cells.Item(31, 15).Value = 17000
$ws.Cells.Item(31, 16).Value = 17000
$ws.Cells.Item(31, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(31, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 19).Value = 1700
$ws.Cells.Item(31, 20).Value = 10
